# The workbook tracks LDLC phone prices with one timestamped snapshot
# column per scrape, followed by "nom" (name) and "url_produit" (url)
# columns. This edit adds a brand-new snapshot column ("2026-01-30
# 06:34:17") right before the "nom"/"url_produit" columns, pushing them
# one column to the right (BD -> BE, BE -> BF), and fills the new
# snapshot column with the latest known price (copied from the former
# last price column, now one column to the left of the new column)
# for every row that had a price, leaving it blank otherwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at BD; this shifts the existing BD (nom) and BE
# (url_produit) columns one position to the right, to BE and BF.
$ws.Columns("BD").Insert()

# New header for the inserted column.
$ws.Range("BD1").Value2 = "2026-01-30 06:34:17"

# The previous last price column is now BC (column 55); the freshly
# inserted column is BD (column 56).
$lastPriceCol = 55
$newCol = 56
$firstDataRow = 2
$lastDataRow = 206

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $price = $ws.Cells.Item($r, $lastPriceCol).Value2
    if ($price -ne $null -and $price -ne "") {
        $ws.Cells.Item($r, $newCol).Value2 = $price
    } else {
        # No price recorded for this product at this snapshot time;
        # still materialize the cell (blank), matching the previously
        # blank price columns for this row.
        $ws.Cells.Item($r, $newCol).Value2 = ""
    }
}
